$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.655.65"
Set-TextValue $ws.Range("D3") "2.122.86"
Set-TextValue $ws.Range("E3") "  +0.69%  "
Set-TextValue $ws.Range("D4") "1.014"
Set-TextValue $ws.Range("E4") "  +0.67%  "
Set-TextValue $ws.Range("D5") "353.12"
Set-TextValue $ws.Range("E5") "  +5.32%  "
Set-TextValue $ws.Range("D6") "1.012"
Set-TextValue $ws.Range("E6") "  +0.47%  "
Set-TextValue $ws.Range("D7") "0.5283"
Set-TextValue $ws.Range("E7") "  +0.99%  "
Set-TextValue $ws.Range("D8") "0.4533"
Set-TextValue $ws.Range("E8") "  -0.35%  "
Set-TextValue $ws.Range("D9") "53.95"
Set-TextValue $ws.Range("E9") "  +1.63%  "
Set-TextValue $ws.Range("D10") "0.09087"
Set-TextValue $ws.Range("E10") "  +1.68%  "
Set-TextValue $ws.Range("D11") "1.182"
Set-TextValue $ws.Range("E11") "  +0.47%  "
Set-TextValue $ws.Range("D12") "24.59"
Set-TextValue $ws.Range("E12") "  +1.08%  "
Set-TextValue $ws.Range("D13") "2.128.86"
Set-TextValue $ws.Range("E13") "  +1.25%  "
Set-TextValue $ws.Range("D14") "6.851"
Set-TextValue $ws.Range("E14") "  -0.20%  "
Set-TextValue $ws.Range("E15") "  +0.49%  "
Set-TextValue $ws.Range("D16") "102.43"
Set-TextValue $ws.Range("E16") "  +6.05%  "
Set-TextValue $ws.Range("E17") "  +2.82%  "
Set-TextValue $ws.Range("E18") "  +0.51%  "
Set-TextValue $ws.Range("D19") "0.06721"
Set-TextValue $ws.Range("E19") "  +0.76%  "
Set-TextValue $ws.Range("D20") "19.44"
Set-TextValue $ws.Range("E20") "  +1.11%  "
Set-TextValue $ws.Range("E21") "  +0.51%  "
Set-TextValue $ws.Range("D22") "6.339"
Set-TextValue $ws.Range("E22") "  -0.35%  "
Set-TextValue $ws.Range("D23") "30.715.41"
Set-TextValue $ws.Range("E23") "  +0.66%  "
Set-TextValue $ws.Range("E24") "  +3.19%  "
Set-TextValue $ws.Range("D25") "2.390"
Set-TextValue $ws.Range("E25") "  +0.73%  "
Set-TextValue $ws.Range("D26") "2.370.47"
Set-TextValue $ws.Range("E26") "  +0.73%  "
Set-TextValue $ws.Range("D27") "22.43"
Set-TextValue $ws.Range("E27") "  +0.28%  "
Set-TextValue $ws.Range("B28") "Monero"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "165.34"
Set-TextValue $ws.Range("E28") "  +0.97%  "
Set-TextValue $ws.Range("B29") "LidoDAOToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.570"
Set-TextValue $ws.Range("E29") "  +0.96%  "
Set-TextValue $ws.Range("D30") "136.41"
Set-TextValue $ws.Range("E30") "  +2.11%  "
Set-TextValue $ws.Range("D31") "1.198"
Set-TextValue $ws.Range("E31") "  -2.18%  "
Set-TextValue $ws.Range("D32") "0.1079"
Set-TextValue $ws.Range("E32") "  +0.56%  "
Set-TextValue $ws.Range("D33") "1.655"
Set-TextValue $ws.Range("E33") "  -0.84%  "
Set-TextValue $ws.Range("D34") "6.369"
Set-TextValue $ws.Range("E34") "  +0.56%  "
Set-TextValue $ws.Range("D35") "4.027"
Set-TextValue $ws.Range("E35") "  +1.98%  "
Set-TextValue $ws.Range("D36") "6.015"
Set-TextValue $ws.Range("E36") "  +5.59%  "
Set-TextValue $ws.Range("D37") "10.34"
Set-TextValue $ws.Range("E37") "  -1.33%  "
Set-TextValue $ws.Range("D38") "0.02660"
Set-TextValue $ws.Range("E38") "  +2.91%  "
Set-TextValue $ws.Range("D39") "0.06889"
Set-TextValue $ws.Range("E39") "  +0.66%  "
Set-TextValue $ws.Range("D40") "0.2321"
Set-TextValue $ws.Range("E40") "  +0.72%  "
Set-TextValue $ws.Range("E41") "  -0.87%  "
Set-TextValue $ws.Range("D42") "0.6928"
Set-TextValue $ws.Range("E42") "  +0.56%  "
Set-TextValue $ws.Range("D43") "1.279"
Set-TextValue $ws.Range("E43") "  +2.44%  "
Set-TextValue $ws.Range("D44") "14.76"
Set-TextValue $ws.Range("E44") "  +5.27%  "
Set-TextValue $ws.Range("E45") "  +1.41%  "
Set-TextValue $ws.Range("D46") "2.334"
Set-TextValue $ws.Range("E46") "  -0.22%  "
Set-TextValue $ws.Range("D47") "3.772"
Set-TextValue $ws.Range("E47") "  +2.66%  "
Set-TextValue $ws.Range("D48") "0.00000000367"
Set-TextValue $ws.Range("E48") "  +9.74%  "
Set-TextValue $ws.Range("D49") "1.258"
Set-TextValue $ws.Range("E49") "  +0.60%  "
Set-TextValue $ws.Range("D50") "82.94"
Set-TextValue $ws.Range("E50") "  -0.59%  "
Set-TextValue $ws.Range("D51") "0.07310"
Set-TextValue $ws.Range("E51") "  +2.35%  "
